$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# Rearranged components in BOM: swap the two component rows (22/23)
# Row22 was U2/MAX3232CSE..., Row23 was U1/STM32G0B1CEU6...
# After the edit, Row22 holds the U1/STM32 data and Row23 holds the
# U2/MAX3232 data.
# ------------------------------------------------------------------
$a22 = $ws.Range("A22").Value2
$b22 = $ws.Range("B22").Value2
$c22 = $ws.Range("C22").Value2
$d22 = $ws.Range("D22").Value2
$e22 = $ws.Range("E22").Value2
$f22 = $ws.Range("F22").Value2

$a23 = $ws.Range("A23").Value2
$b23 = $ws.Range("B23").Value2
$c23 = $ws.Range("C23").Value2
$d23 = $ws.Range("D23").Value2
$e23 = $ws.Range("E23").Value2
$f23 = $ws.Range("F23").Value2

$ws.Range("A22").Value2 = $a23
$ws.Range("B22").Value2 = $b23
$ws.Range("C22").Value2 = $c23
$ws.Range("D22").Value2 = $d23
$ws.Range("E22").Value2 = $e23
$ws.Range("F22").Value2 = $f23

$ws.Range("A23").Value2 = $a22
$ws.Range("B23").Value2 = $b22
$ws.Range("C23").Value2 = $c22
$ws.Range("D23").Value2 = $d22
$ws.Range("E23").Value2 = $e22
$ws.Range("F23").Value2 = $f22

# ------------------------------------------------------------------
# The hyperlinks on E22/E23 need to move together with their data,
# i.e. the MAX3232 hyperlink goes to E23 and the STM32 hyperlink
# goes to E22 (they stay keyed to the part that now sits in that row).
# ------------------------------------------------------------------
$urlE22_old = "https://octopart.com/max3232cse%2B-analog+devices-124215225"
$urlE23_old = "https://octopart.com/stm32g0b1ceu6-stmicroelectronics-116363364?r=sp"

$ws.Hyperlinks.Add($ws.Range("E22"), $urlE23_old) | Out-Null
$ws.Hyperlinks.Add($ws.Range("E23"), $urlE22_old) | Out-Null

# ------------------------------------------------------------------
# Row 21 no longer needs an explicit (wrapped-text) row height -
# let Excel auto fit it back to the default height.
# ------------------------------------------------------------------
$ws.Rows.Item(21).AutoFit()

# ------------------------------------------------------------------
# Update the view: scrolled down a bit and row 22 is fully selected.
# ------------------------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("A22:XFD22").Select() | Out-Null

Write-Host "done"
